$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "TODO Before 0.0.1": mark row 17/18 item as done, add AutoFilter that
# shows only "in-progress"/"todo" rows (hides "done"/"done " rows), and move
# the selection to B30.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TODO Before 0.0.1")

# C18 was blank; fill in its status (copy the format from the cell above so
# it keeps the same centered style used by the rest of column C).
$ws1.Range("C17").Copy()
$ws1.Range("C18").PasteSpecial(-4122)
$ws1.Range("C18").Value = "done"

# Turn on the AutoFilter for the table, filtering column C ("status") to
# only show "in-progress" and "todo" rows - this also hides the "done" rows.
$ws1.Range("A1:E46").AutoFilter(3, @("in-progress", "todo"), 7)

# Register the (hidden) _FilterDatabase defined name that Excel creates for
# the sheet once an AutoFilter is applied.
$fd1 = $ws1.Names.Add("_xlnm._FilterDatabase", "='TODO Before 0.0.1'!`$A`$1:`$E`$46")
$fd1.Visible = $false

$ws1.Activate()
$ws1.Range("B30").Select()

# ---------------------------------------------------------------------------
# Sheet "Logs": fix up the date on the last existing entry, append the two
# new log entries, and move the selection to B59.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Logs")

# The last existing entry's date was corrected (moved from 2024-10-01 to
# 2024-10-10).
$ws3.Range("A56").Value = 45575

$ws3.Range("A56").Copy()
$ws3.Range("A57").PasteSpecial(-4122)
$ws3.Range("A57").Value = 45576
$ws3.Range("B57").Value = "tried to fix pickup items handling, fails. Back from vacation :) "

$ws3.Range("A56").Copy()
$ws3.Range("A58").PasteSpecial(-4122)
$ws3.Range("A58").Value = 45577
$ws3.Range("B58").Value = "start work on item info block - skills need to be clasified  (notes in notebook)"

$ws3.Activate()
$ws3.Range("B59").Select()
